$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sender ("from") block
$ws.Range("B2").Value = "Wilkinson, Wilkinson and Wilkinson"
$ws.Range("B3").Value = "%6961 Britney Meadow`nLake Willieville, SC 83113-3279, reilly.everett@gmail.com"
$ws.Range("B4").Value = "https://uny.com/gjly-gnet.html"
$ws.Range("B5").Value = "{0x140004040c0}"

# Bill To / Ship To block
$ws.Range("B10").Value = "Mr. Sofia Metz"
$ws.Range("D10").Value = "Murray Marquardt II"

$ws.Range("B11").Value = "Ziemann and Sons"
$ws.Range("D11").Value = "Ziemann and Sons"

$ws.Range("B12").Value = "%71 Sadye Light`nFletatown, KY 53478"
$ws.Range("D12").Value = "%71 Sadye Light`nFletatown, KY 53478"

$ws.Range("B13").Value = "dillan.nader@rnf.com"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "+18002574390"
